$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Recorded By" column (G) lists System, dnasr281@gmail.com in some
# unspecified order; swap the order to dnasr281@gmail.com, System for every
# row where it currently reads "System, dnasr281@gmail.com".
$target = "System, dnasr281@gmail.com"
$replacement = "dnasr281@gmail.com, System"

$lastRow = $ws.UsedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Range("G$row")
    if ($cell.Value2 -eq $target) {
        $cell.Value2 = $replacement
    }
}
